$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the stray "_GoBack" bookmark that wraps the title paragraph
#    ("PHP EXERCÍCIOS - CICLOS"). Bookmark.Delete() is a no-op in this
#    runtime, so we rebuild the paragraph (same text/formatting, no
#    bookmark) via Range.InsertXML, which replaces both the
#    <w:bookmarkStart> inside it and the <w:bookmarkEnd> that immediately
#    follows it (still inside the replaced range, since it sits right
#    before the next paragraph).
# ---------------------------------------------------------------------------

$titleText = "PHP EXERCÍCIOS - CICLOS"
$titlePara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs($i)
    if ($cand.Range.Text.TrimEnd([char]13, [char]7) -eq $titleText) {
        $titlePara = $cand
        break
    }
}

if ($titlePara -ne $null) {
    # Range covering the title paragraph *and* the following paragraph mark,
    # so the loose <w:bookmarkEnd/> that sits between the two <w:p> elements
    # is included in what gets replaced.
    $titleRange = $d.Range($titlePara.Range.Start, $titlePara.Range.End + 1)

    $titleXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="left"/><w:rPr><w:rFonts w:hint="default"/><w:b/><w:bCs/><w:sz w:val="36"/><w:szCs w:val="36"/><w:lang w:val="pt-PT"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="default"/><w:b/><w:bCs/><w:sz w:val="36"/><w:szCs w:val="36"/><w:lang w:val="pt-PT"/></w:rPr><w:t>PHP EXERCÍCIOS - CICLOS</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="left"/><w:rPr><w:rFonts w:hint="default"/><w:b/><w:bCs/><w:sz w:val="36"/><w:szCs w:val="36"/><w:lang w:val="pt-PT"/></w:rPr></w:pPr></w:p>'

    $titleRange.InsertXML($titleXml)
}

Write-Host "GoBack still on title?" $d.Bookmarks.Exists("_GoBack")

# ---------------------------------------------------------------------------
# 2) Split the "factorial" paragraph's run into two runs (after "o fa", the
#    "c" is lost) and put a fresh, empty "_GoBack" bookmark exactly at the
#    split point - reproducing Word's "last edit position" marker.
# ---------------------------------------------------------------------------

$oldText = "4 - Escreva um programa que calcule o factorial de um número."
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs($i)
    if ($cand.Range.Text.TrimEnd([char]13, [char]7) -eq $oldText) {
        $targetPara = $cand
        break
    }
}

if ($targetPara -ne $null) {
    $factorialRange = $d.Range($targetPara.Range.Start, $targetPara.Range.End)

    $factorialXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="4"/><w:keepNext w:val="0"/><w:keepLines w:val="0"/><w:widowControl/><w:suppressLineNumbers w:val="0"/><w:rPr><w:rFonts w:hint="default" w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:cstheme="minorBidi"/><w:kern w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en-US" w:eastAsia="zh-CN" w:bidi="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="default" w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:cstheme="minorBidi"/><w:kern w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en-US" w:eastAsia="zh-CN" w:bidi="ar-SA"/></w:rPr><w:t>4 - Escreva um programa que calcule o fa</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:rFonts w:hint="default" w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:cstheme="minorBidi"/><w:kern w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en-US" w:eastAsia="zh-CN" w:bidi="ar-SA"/></w:rPr><w:t>torial de um número.</w:t></w:r></w:p>'

    $factorialRange.InsertXML($factorialXml)
}

Write-Host "GoBack present now?" $d.Bookmarks.Exists("_GoBack")
$gb = $d.Bookmarks.Item("_GoBack")
Write-Host "GoBack range:" $gb.Start $gb.End
